# Update the "TestOverview" sheet (Tabelle1):
#  - add two new test rows (LoadBalancing, UniReceive) to the results table
#  - mark several existing/new rows' multiplicity-check columns as "correct"
#  - widen the new Comment-adjacent column and fix up the view/selection

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")
$ws.Activate()

# --- existing rows: flag model-checking columns E/F as "correct" ---
$ws.Range("E15").Value = "correct"
$ws.Range("F15").Value = "correct"

$ws.Range("E20").Value = "correct"
$ws.Range("F20").Value = "correct"

# --- new row 22: LoadBalancing ---
$ws.Range("A22").Value = "LoadBalancing"
$ws.Range("B22").Value = "no"
$ws.Range("C22").Value = "no"
$ws.Range("D22").Value = "no"
$ws.Range("E22").Value = "correct"
$ws.Range("F22").Value = "correct"

# --- new row 23: UniReceive ---
$ws.Range("A23").Value = "UniReceive"
$ws.Range("B23").Value = "no"
$ws.Range("C23").Value = "no"
$ws.Range("D23").Value = "no"

# --- column G width (bestFit-style width for the new comment column) ---
$ws.Columns.Item(7).ColumnWidth = 10.43

# --- refresh the view: scroll & selection ---
$excel.ActiveWindow.ScrollRow = 7
$ws.Range("A22").Select()
